$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ntf3"
$ws.Cells.Item(2,3).Value = "Ntrk3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.679513333333333
$ws.Cells.Item(2,8).Value = 5.038539999999999
$ws.Cells.Item(2,9).Value = 0.376631045782902
$ws.Cells.Item(2,10).Value = 0.3928053077137587
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.002145666666666667
$ws.Cells.Item(2,14).Value = 0.006437
$ws.Cells.Item(2,15).Value = 0.00807320947388686
$ws.Cells.Item(2,16).Value = 0.008843624333499573
$ws.Cells.Item(2,17).Value = 0.003603675775555555
$ws.Cells.Item(2,18).Value = 0.03243308197999999
$ws.Cells.Item(2,19).Value = 0.00304062132697444
$ws.Cells.Item(2,20).Value = 0.003473822577625184

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ntf3"
$ws.Cells.Item(3,3).Value = "Ntrk3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.679513333333333
$ws.Cells.Item(3,8).Value = 5.038539999999999
$ws.Cells.Item(3,9).Value = 0.376631045782902
$ws.Cells.Item(3,10).Value = 0.3928053077137587
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.03558433333333334
$ws.Cells.Item(3,14).Value = 0.106753
$ws.Cells.Item(3,15).Value = 0.1338883534202026
$ws.Cells.Item(3,16).Value = 0.1466651279282398
$ws.Cells.Item(3,17).Value = 0.05976436229111111
$ws.Cells.Item(3,18).Value = 0.5378792606199999
$ws.Cells.Item(3,19).Value = 0.05042651056680168
$ws.Cells.Item(3,20).Value = 0.05761084070673003

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ntf3"
$ws.Cells.Item(4,3).Value = "Ntrk3"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.679513333333333
$ws.Cells.Item(4,8).Value = 5.038539999999999
$ws.Cells.Item(4,9).Value = 0.376631045782902
$ws.Cells.Item(4,10).Value = 0.3928053077137587
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.04573
$ws.Cells.Item(4,14).Value = 0.13719
$ws.Cells.Item(4,15).Value = 0.172062079807758
$ws.Cells.Item(4,16).Value = 0.188481718551003
$ws.Cells.Item(4,17).Value = 0.07680414473333333
$ws.Cells.Item(4,18).Value = 0.6912373025999999
$ws.Cells.Item(4,19).Value = 0.06480392105757704
$ws.Cells.Item(4,20).Value = 0.07403661945384481

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Ntf3"
$ws.Cells.Item(5,3).Value = "Ntrk3"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.679513333333333
$ws.Cells.Item(5,8).Value = 5.038539999999999
$ws.Cells.Item(5,9).Value = 0.376631045782902
$ws.Cells.Item(5,10).Value = 0.3928053077137587
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.06945950000000001
$ws.Cells.Item(5,14).Value = 0.138919
$ws.Cells.Item(5,15).Value = 0.2613458568206203
$ws.Cells.Item(5,16).Value = 0.1908571459974254
$ws.Cells.Item(5,17).Value = 0.1166581563766667
$ws.Cells.Item(5,18).Value = 0.6999489382599999
$ws.Cells.Item(5,19).Value = 0.09843096336537879
$ws.Cells.Item(5,20).Value = 0.07496969996288845

$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Ntf3"
$ws.Cells.Item(6,3).Value = "Ntrk3"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.679513333333333
$ws.Cells.Item(6,8).Value = 5.038539999999999
$ws.Cells.Item(6,9).Value = 0.376631045782902
$ws.Cells.Item(6,10).Value = 0.3928053077137587
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.1128566666666667
$ws.Cells.Item(6,14).Value = 0.33857
$ws.Cells.Item(6,15).Value = 0.4246305004775321
$ws.Cells.Item(6,16).Value = 0.4651523831898322
$ws.Cells.Item(6,17).Value = 0.1895442764222222
$ws.Cells.Item(6,18).Value = 1.7058984878
$ws.Cells.Item(6,19).Value = 0.15992902946617
$ws.Cells.Item(6,20).Value = 0.1827143250126702

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ntf3"
$ws.Cells.Item(7,3).Value = "Ntrk3"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.228940333333334
$ws.Cells.Item(7,8).Value = 6.686821
$ws.Cells.Item(7,9).Value = 0.4998401096732527
$ws.Cells.Item(7,10).Value = 0.5213055330575571
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.002145666666666667
$ws.Cells.Item(7,14).Value = 0.006437
$ws.Cells.Item(7,15).Value = 0.00807320947388686
$ws.Cells.Item(7,16).Value = 0.008843624333499573
$ws.Cells.Item(7,17).Value = 0.004782562975222223
$ws.Cells.Item(7,18).Value = 0.043043066777
$ws.Cells.Item(7,19).Value = 0.004035313908842751
$ws.Cells.Item(7,20).Value = 0.004610230297335778

$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Ntf3"
$ws.Cells.Item(8,3).Value = "Ntrk3"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.228940333333334
$ws.Cells.Item(8,8).Value = 6.686821
$ws.Cells.Item(8,9).Value = 0.4998401096732527
$ws.Cells.Item(8,10).Value = 0.5213055330575571
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.03558433333333334
$ws.Cells.Item(8,14).Value = 0.106753
$ws.Cells.Item(8,15).Value = 0.1338883534202026
$ws.Cells.Item(8,16).Value = 0.1466651279282398
$ws.Cells.Item(8,17).Value = 0.07931535580144446
$ws.Cells.Item(8,18).Value = 0.713838202213
$ws.Cells.Item(8,19).Value = 0.06692276925752527
$ws.Cells.Item(8,20).Value = 0.07645734269558586

$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Ntf3"
$ws.Cells.Item(9,3).Value = "Ntrk3"
$ws.Cells.Item(9,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.228940333333334
$ws.Cells.Item(9,8).Value = 6.686821
$ws.Cells.Item(9,9).Value = 0.4998401096732527
$ws.Cells.Item(9,10).Value = 0.5213055330575571
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.04573
$ws.Cells.Item(9,14).Value = 0.13719
$ws.Cells.Item(9,15).Value = 0.172062079807758
$ws.Cells.Item(9,16).Value = 0.188481718551003
$ws.Cells.Item(9,17).Value = 0.1019294414433333
$ws.Cells.Item(9,18).Value = 0.9173649729900001
$ws.Cells.Item(9,19).Value = 0.08600352884171772
$ws.Cells.Item(9,20).Value = 0.09825656276083508

$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Ntf3"
$ws.Cells.Item(10,3).Value = "Ntrk3"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 2.228940333333334
$ws.Cells.Item(10,8).Value = 6.686821
$ws.Cells.Item(10,9).Value = 0.4998401096732527
$ws.Cells.Item(10,10).Value = 0.5213055330575571
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.06945950000000001
$ws.Cells.Item(10,14).Value = 0.138919
$ws.Cells.Item(10,15).Value = 0.2613458568206203
$ws.Cells.Item(10,16).Value = 0.1908571459974254
$ws.Cells.Item(10,17).Value = 0.1548210810831667
$ws.Cells.Item(10,18).Value = 0.9289264864990001
$ws.Cells.Item(10,19).Value = 0.1306311417358691
$ws.Cells.Item(10,20).Value = 0.09949488623203183

$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Ntf3"
$ws.Cells.Item(11,3).Value = "Ntrk3"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 2.228940333333334
$ws.Cells.Item(11,8).Value = 6.686821
$ws.Cells.Item(11,9).Value = 0.4998401096732527
$ws.Cells.Item(11,10).Value = 0.5213055330575571
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.1128566666666667
$ws.Cells.Item(11,14).Value = 0.33857
$ws.Cells.Item(11,15).Value = 0.4246305004775321
$ws.Cells.Item(11,16).Value = 0.4651523831898322
$ws.Cells.Item(11,17).Value = 0.2515507762188889
$ws.Cells.Item(11,18).Value = 2.26395698597
$ws.Cells.Item(11,19).Value = 0.2122473559292978
$ws.Cells.Item(11,20).Value = 0.2424865110717686

$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Ntf3"
$ws.Cells.Item(12,3).Value = "Ntrk3"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.550853
$ws.Cells.Item(12,8).Value = 1.101706
$ws.Cells.Item(12,9).Value = 0.1235288445438454
$ws.Cells.Item(12,10).Value = 0.08588915922868416
$ws.Cells.Item(12,11).Value = 1
$ws.Cells.Item(12,12).Value = 0.3333333333333333
$ws.Cells.Item(12,13).Value = 0.002145666666666667
$ws.Cells.Item(12,14).Value = 0.006437
$ws.Cells.Item(12,15).Value = 0.00807320947388686
$ws.Cells.Item(12,16).Value = 0.008843624333499573
$ws.Cells.Item(12,17).Value = 0.001181946920333334
$ws.Cells.Item(12,18).Value = 0.007091681522
$ws.Cells.Item(12,19).Value = 0.0009972742380696699
$ws.Cells.Item(12,20).Value = 0.0007595714585386106

$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Ntf3"
$ws.Cells.Item(13,3).Value = "Ntrk3"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.550853
$ws.Cells.Item(13,8).Value = 1.101706
$ws.Cells.Item(13,9).Value = 0.1235288445438454
$ws.Cells.Item(13,10).Value = 0.08588915922868416
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.03558433333333334
$ws.Cells.Item(13,14).Value = 0.106753
$ws.Cells.Item(13,15).Value = 0.1338883534202026
$ws.Cells.Item(13,16).Value = 0.1466651279282398
$ws.Cells.Item(13,17).Value = 0.01960173676966667
$ws.Cells.Item(13,18).Value = 0.117610420618
$ws.Cells.Item(13,19).Value = 0.01653907359587563
$ws.Cells.Item(13,20).Value = 0.01259694452592392

$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,2).Value = "Ntf3"
$ws.Cells.Item(14,3).Value = "Ntrk3"
$ws.Cells.Item(14,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.550853
$ws.Cells.Item(14,8).Value = 1.101706
$ws.Cells.Item(14,9).Value = 0.1235288445438454
$ws.Cells.Item(14,10).Value = 0.08588915922868416
$ws.Cells.Item(14,11).Value = 1
$ws.Cells.Item(14,12).Value = 0.3333333333333333
$ws.Cells.Item(14,13).Value = 0.04573
$ws.Cells.Item(14,14).Value = 0.13719
$ws.Cells.Item(14,15).Value = 0.172062079807758
$ws.Cells.Item(14,16).Value = 0.188481718551003
$ws.Cells.Item(14,17).Value = 0.02519050769
$ws.Cells.Item(14,18).Value = 0.15114304614
$ws.Cells.Item(14,19).Value = 0.02125462990846326
$ws.Cells.Item(14,20).Value = 0.01618853633632313

$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,2).Value = "Ntf3"
$ws.Cells.Item(15,3).Value = "Ntrk3"
$ws.Cells.Item(15,4).Value = "MuSCs"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.550853
$ws.Cells.Item(15,8).Value = 1.101706
$ws.Cells.Item(15,9).Value = 0.1235288445438454
$ws.Cells.Item(15,10).Value = 0.08588915922868416
$ws.Cells.Item(15,11).Value = 2
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.06945950000000001
$ws.Cells.Item(15,14).Value = 0.138919
$ws.Cells.Item(15,15).Value = 0.2613458568206203
$ws.Cells.Item(15,16).Value = 0.1908571459974254
$ws.Cells.Item(15,17).Value = 0.03826197395350001
$ws.Cells.Item(15,18).Value = 0.153047895814
$ws.Cells.Item(15,19).Value = 0.03228375171937248
$ws.Cells.Item(15,20).Value = 0.01639255980250509

$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,2).Value = "Ntf3"
$ws.Cells.Item(16,3).Value = "Ntrk3"
$ws.Cells.Item(16,4).Value = "Resolving-Mac"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.550853
$ws.Cells.Item(16,8).Value = 1.101706
$ws.Cells.Item(16,9).Value = 0.1235288445438454
$ws.Cells.Item(16,10).Value = 0.08588915922868416
$ws.Cells.Item(16,11).Value = 1
$ws.Cells.Item(16,12).Value = 0.3333333333333333
$ws.Cells.Item(16,13).Value = 0.1128566666666667
$ws.Cells.Item(16,14).Value = 0.33857
$ws.Cells.Item(16,15).Value = 0.4246305004775321
$ws.Cells.Item(16,16).Value = 0.4651523831898322
$ws.Cells.Item(16,17).Value = 0.06216743340333333
$ws.Cells.Item(16,18).Value = 0.37300460042
$ws.Cells.Item(16,19).Value = 0.05245411508206433
$ws.Cells.Item(16,20).Value = 0.03995154710539341
